{"js": "// Update the worksheet date and every three-digit x one-digit multiplication\n// fact in the table to the new day's values.\nconst replacements = [\n  [\"2024-08-25 Sunday\", \"2024-08-26 Monday\"],\n  [\"539\u00d79=4851\", \"668\u00d76=4008\"],\n  [\"462\u00d73=1386\", \"910\u00d72=1820\"],\n  [\"415\u00d75=2075\", \"227\u00d74=908\"],\n  [\"947\u00d77=6629\", \"242\u00d78=1936\"],\n  [\"246\u00d75=1230\", \"232\u00d74=928\"],\n  [\"321\u00d72=642\", \"217\u00d74=868\"],\n  [\"246\u00d78=1968\", \"993\u00d76=5958\"],\n  [\"782\u00d77=5474\", \"468\u00d75=2340\"],\n  [\"984\u00d78=7872\", \"751\u00d75=3755\"],\n  [\"887\u00d76=5322\", \"664\u00d72=1328\"],\n  [\"267\u00d72=534\", \"495\u00d72=990\"],\n  [\"824\u00d78=6592\", \"769\u00d78=6152\"],\n  [\"897\u00d78=7176\", \"675\u00d76=4050\"],\n  [\"547\u00d78=4376\", \"640\u00d75=3200\"],\n  [\"630\u00d74=2520\", \"994\u00d78=7952\"],\n  [\"731\u00d74=2924\", \"161\u00d74=644\"],\n  [\"981\u00d74=3924\", \"669\u00d73=2007\"],\n  [\"682\u00d76=4092\", \"433\u00d79=3897\"],\n  [\"210\u00d75=1050\", \"224\u00d79=2016\"],\n  [\"249\u00d76=1494\", \"295\u00d79=2655\"],\n  [\"839\u00d75=4195\", \"708\u00d76=4248\"],\n  [\"413\u00d77=2891\", \"972\u00d79=8748\"],\n  [\"209\u00d73=627\", \"207\u00d77=1449\"],\n  [\"858\u00d76=5148\", \"826\u00d75=4130\"],\n  [\"714\u00d77=4998\", \"636\u00d76=3816\"],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every three-digit x one-digit multiplication\n# fact in the table to the new day's values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-08-25 Sunday\", \"2024-08-26 Monday\")\n    ,@(\"539\u00d79=4851\", \"668\u00d76=4008\")\n    ,@(\"462\u00d73=1386\", \"910\u00d72=1820\")\n    ,@(\"415\u00d75=2075\", \"227\u00d74=908\")\n    ,@(\"947\u00d77=6629\", \"242\u00d78=1936\")\n    ,@(\"246\u00d75=1230\", \"232\u00d74=928\")\n    ,@(\"321\u00d72=642\", \"217\u00d74=868\")\n    ,@(\"246\u00d78=1968\", \"993\u00d76=5958\")\n    ,@(\"782\u00d77=5474\", \"468\u00d75=2340\")\n    ,@(\"984\u00d78=7872\", \"751\u00d75=3755\")\n    ,@(\"887\u00d76=5322\", \"664\u00d72=1328\")\n    ,@(\"267\u00d72=534\", \"495\u00d72=990\")\n    ,@(\"824\u00d78=6592\", \"769\u00d78=6152\")\n    ,@(\"897\u00d78=7176\", \"675\u00d76=4050\")\n    ,@(\"547\u00d78=4376\", \"640\u00d75=3200\")\n    ,@(\"630\u00d74=2520\", \"994\u00d78=7952\")\n    ,@(\"731\u00d74=2924\", \"161\u00d74=644\")\n    ,@(\"981\u00d74=3924\", \"669\u00d73=2007\")\n    ,@(\"682\u00d76=4092\", \"433\u00d79=3897\")\n    ,@(\"210\u00d75=1050\", \"224\u00d79=2016\")\n    ,@(\"249\u00d76=1494\", \"295\u00d79=2655\")\n    ,@(\"839\u00d75=4195\", \"708\u00d76=4248\")\n    ,@(\"413\u00d77=2891\", \"972\u00d79=8748\")\n    ,@(\"209\u00d73=627\", \"207\u00d77=1449\")\n    ,@(\"858\u00d76=5148\", \"826\u00d75=4130\")\n    ,@(\"714\u00d77=4998\", \"636\u00d76=3816\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
